# Add the AU logo picture (already embedded in the deck, used on slide 1)
# to the last slide ("S.O.L.I.D principles"), matching the commit
# "Added AU_LOGO to Slides".
#
# The logo is cloned from the existing picture shape on slide 1 (instead of
# Shapes.AddPicture against a file) so the already-embedded image bytes /
# relationship get reused exactly like PowerPoint does for a copy-paste of
# an existing picture, without needing filesystem access.

$p = $ppt.ActivePresentation

# Slide that already carries the AU logo picture (slide 1, "Billede 3" - a
# msoPicture shape, Type 13).
$srcSlide = $p.Slides.Item(1)
$logoShape = $null
for ($i = 1; $i -le $srcSlide.Shapes.Count; $i++) {
    $candidate = $srcSlide.Shapes.Item($i)
    if ($candidate.Type -eq 13) {
        $logoShape = $candidate
        break
    }
}

# Destination slide: the last slide in the deck ("S.O.L.I.D principles").
$dstSlide = $p.Slides.Item($p.Slides.Count)

$logoShape.Copy()
$pasted = $dstSlide.Shapes.Paste()
$pic = $pasted.Item(1)

$pic.Name = "Billede 4"

# Position/size (target EMU taken from the edited slide) converted to points
# since the PowerPoint object model works in points (1 pt = 12700 EMU).
# Shape.Left/Top/Width/Height are stored as single-precision floats, so add
# half an EMU (in points) before PowerPoint truncates back to EMU on save -
# this rounds to the nearest EMU instead of always truncating down.
$emuPerPt = 12700.0
$halfEmuPt = 0.5 / $emuPerPt

$pic.Left = (9077642 / $emuPerPt) + $halfEmuPt
$pic.Top = (5408502 / $emuPerPt) + $halfEmuPt
$pic.Width = (2276158 / $emuPerPt) + $halfEmuPt
$pic.Height = (947848 / $emuPerPt) + $halfEmuPt
